# Power system workbook edit
# Reproduces: new "Battery Standby Voltage" row in Battery Monitoring,
# helper Vout columns, new Resolution(mV) cell, Power Budgeting runtime
# calc + removal of the Motors "Nominal A" value, and cosmetic
# selection / active-tab changes on the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Battery Monitoring"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Battery Monitoring")

# Insert a new row 3 ("Battery Standby Voltage") above the existing
# "Battery Nominal Voltage" row, copying the row-2 layout/format down
# so the new row matches the rest of the voltage table, then shift the
# old row 3 (and everything below it) down automatically.
$ws1.Rows("3:3").Insert() | Out-Null
$ws1.Range("A2:E2").Copy() | Out-Null
$ws1.Range("A3:E3").PasteSpecial(-4122) | Out-Null

$ws1.Range("A3").Value = "Battery Standby Voltage"
$ws1.Range("B3").Value = 13.8
$ws1.Range("C3").Formula = "=(B3*`$B`$12)/(`$B`$11+`$B`$12)"
$ws1.Range("D3").Formula = "=(1/`$B`$10)*C3"
$ws1.Range("E3").Formula = "=D3-D5"

# New helper "Vout" column next to the voltage table.
$ws1.Range("H2").Value = "Vout"
$ws1.Range("H4").Value = 1.738
$ws1.Range("I4").Formula = "=(B11+B12)*H4/B12"

# New "Resolution (V/Point) x1000" helper cell below the resolution row
# (now row 10 after the insert above).
$ws1.Range("C10").Formula = "=B10*1000"
$ws1.Range("C10").NumberFormat = "0.0000"

# Reposition the circuit picture so it still starts/ends at the same
# rows relative to the table after the row insert above (was rows
# 19-38 1-based, now 20-39 1-based / 18-37 -> 19-38 zero based).
$shp1 = $ws1.Shapes.Item(1)
$shp1.Top = $shp1.Top + 15

# Make "Battery Monitoring" the active sheet/selection.
$ws1.Activate()
$ws1.Range("G6").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet 2: "Power Budgeting"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Power Budgeting")

# The "Motors" nominal-amp figure was cleared out (so the dependent
# totals drop to 0 for that row).
$ws2.Range("B10").ClearContents() | Out-Null

# New runtime estimate row.
$ws2.Range("F15").Value = "Runtime(hr)"
$ws2.Range("F16").Formula = "=20/G16"

$ws2.Columns("B:B").ColumnWidth = 10.28515625
$ws2.Columns("F:F").ColumnWidth = 12.85546875

$ws2.Range("F16").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet 3: "Pins" - selection only changed.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Pins")
$ws3.Range("E19").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet 4: "i2c addr" - selection only changed.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("i2c addr")
$ws4.Range("D22").Select() | Out-Null

# Leave "Battery Monitoring" focused/active, matching the saved file.
$ws1.Activate()
$ws1.Range("G6").Select() | Out-Null
